$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 147.125
$ws.Range("I33").Value = 203.4
$ws.Range("K33").Value = 203.4
$ws.Range("M33").Value = 25.59999999999999

$ws.Range("H132").Value = 3724
$ws.Range("I132").Value = 3537.5557
$ws.Range("J132").Value = 4059.6
$ws.Range("K132").Value = 10612.6671
$ws.Range("L132").Value = 12178.8
$ws.Range("M132").Value = -8082.667099999999
$ws.Range("N132").Value = -17238.8

$ws.Range("H137").Value = 2161.6667
$ws.Range("I137").Value = 2040
$ws.Range("J137").Value = 2283.3333
$ws.Range("K137").Value = 6120
$ws.Range("L137").Value = 6849.999899999999
$ws.Range("M137").Value = -3570
$ws.Range("N137").Value = -11949.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8651.1875
$ws.Range("I32").Value = 4172.7856
$ws.Range("K32").Value = 4172.7856
$ws.Range("M32").Value = -3885.7856

$ws.Range("H61").Value = 2912.5
$ws.Range("I61").Value = 2375
$ws.Range("J61").Value = 3450
$ws.Range("K61").Value = 2375
$ws.Range("L61").Value = 3450
$ws.Range("M61").Value = -2163
$ws.Range("N61").Value = -3874

$ws.Range("H97").Value = 860.13336
$ws.Range("J97").Value = 2666.6667
$ws.Range("L97").Value = 2666.6667
$ws.Range("N97").Value = -3658.6667

$ws.Range("H102").Value = 3196.6
$ws.Range("J102").Value = 3249.5
$ws.Range("L102").Value = 3249.5
$ws.Range("N102").Value = -6493.5

$ws.Range("H136").Value = 2912.5
$ws.Range("I136").Value = 2375
$ws.Range("J136").Value = 3450
$ws.Range("K136").Value = 7125
$ws.Range("L136").Value = 10350
$ws.Range("M136").Value = -4575
$ws.Range("N136").Value = -15450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 10214.143
$ws.Range("I86").Value = 9099.799999999999
$ws.Range("K86").Value = 9099.799999999999
$ws.Range("M86").Value = -7976.799999999999

$ws.Range("H89").Value = 10214.143
$ws.Range("I89").Value = 9099.799999999999
$ws.Range("K89").Value = 45499
$ws.Range("M89").Value = -39883

$ws.Range("H99").Value = 3549.75
$ws.Range("I99").Value = 3549.75
$ws.Range("K99").Value = 3549.75
$ws.Range("M99").Value = -2051.75

$ws.Range("H105").Value = 3763.818
$ws.Range("I105").Value = 3269.125
$ws.Range("J105").Value = 5083
$ws.Range("K105").Value = 3269.125
$ws.Range("L105").Value = 5083
$ws.Range("M105").Value = -1522.125
$ws.Range("N105").Value = -8577

$ws.Range("H107").Value = 888.3333
$ws.Range("J107").Value = 887
$ws.Range("L107").Value = 887
$ws.Range("N107").Value = -4727

$ws.Range("H134").Value = 1430.75
$ws.Range("I134").Value = 1430.75
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4292.25
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -1757.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 643.3333
$ws.Range("I22").Value = 540
$ws.Range("J22").Value = 850
$ws.Range("K22").Value = 540
$ws.Range("L22").Value = 850
$ws.Range("M22").Value = -190
$ws.Range("N22").Value = -1550

$ws.Range("H69").Value = 19500
$ws.Range("I69").Value = 19500
$ws.Range("K69").Value = 19500
$ws.Range("M69").Value = -18751

$ws.Range("H72").Value = 19500
$ws.Range("I72").Value = 19500
$ws.Range("K72").Value = 58500
$ws.Range("M72").Value = -54756

$ws.Range("H105").Value = 1221.75
$ws.Range("I105").Value = 1129.6666
$ws.Range("J105").Value = 1498
$ws.Range("K105").Value = 1129.6666
$ws.Range("L105").Value = 1498
$ws.Range("M105").Value = 617.3334
$ws.Range("N105").Value = -4992

$ws.Range("H116").Value = 16998.5
$ws.Range("J116").Value = 16998.5
$ws.Range("L116").Value = 16998.5
$ws.Range("N116").Value = -26176.5

$ws.Range("H122").Value = 2329
$ws.Range("I122").Value = 2101.6667
$ws.Range("K122").Value = 6305.000100000001
$ws.Range("M122").Value = -3855.000100000001

$ws.Range("H132").Value = 5937
$ws.Range("I132").Value = 6323.5
$ws.Range("J132").Value = 4777.5
$ws.Range("K132").Value = 18970.5
$ws.Range("L132").Value = 14332.5
$ws.Range("M132").Value = -16440.5
$ws.Range("N132").Value = -19392.5

$ws.Range("H134").Value = 1823
$ws.Range("I134").Value = 1823
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5469
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -2934

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1155.2222
$ws.Range("J12").Value = 1484.2858
$ws.Range("L12").Value = 4452.857400000001
$ws.Range("N12").Value = -4798.857400000001

$ws.Range("H50").Value = 195.375
$ws.Range("I50").Value = 160.43478
$ws.Range("K50").Value = 481.30434
$ws.Range("M50").Value = -0.3043399999999679

$ws.Range("H53").Value = 195.375
$ws.Range("I53").Value = 160.43478
$ws.Range("K53").Value = 481.30434
$ws.Range("M53").Value = -0.3043399999999679

$ws.Range("H131").Value = 912.05884

$ws.Range("H140").Value = 1519.75
$ws.Range("I140").Value = 1308.4286
$ws.Range("J140").Value = 2999
$ws.Range("K140").Value = 3925.2858
$ws.Range("L140").Value = 8997
$ws.Range("M140").Value = 1254.7142
$ws.Range("N140").Value = -19357

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 805.1875
$ws.Range("J97").Value = 1011.8571
$ws.Range("L97").Value = 1011.8571
$ws.Range("N97").Value = -2003.8571

$ws.Range("H102").Value = 2200
$ws.Range("I102").Value = 1350
$ws.Range("K102").Value = 1350
$ws.Range("M102").Value = 272

$ws.Range("H132").Value = 3171.35
$ws.Range("I132").Value = 2884.889
$ws.Range("J132").Value = 5749.5
$ws.Range("K132").Value = 8654.667000000001
$ws.Range("L132").Value = 17248.5
$ws.Range("M132").Value = -6124.667000000001
$ws.Range("N132").Value = -22308.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1756.9
$ws.Range("I16").Value = 1807.6666
$ws.Range("J16").Value = 1300
$ws.Range("K16").Value = 1807.6666
$ws.Range("L16").Value = 1300
$ws.Range("M16").Value = -1637.6666
$ws.Range("N16").Value = -1640

$ws.Range("H22").Value = 2050
$ws.Range("I22").Value = 2075
$ws.Range("K22").Value = 2075
$ws.Range("M22").Value = -1780

$ws.Range("H27").Value = 2050
$ws.Range("I27").Value = 2075
$ws.Range("K27").Value = 2075
$ws.Range("M27").Value = -1968

$ws.Range("H40").Value = 6099.75
$ws.Range("I40").Value = 5804.4165
$ws.Range("K40").Value = 5804.4165
$ws.Range("M40").Value = -5668.4165

$ws.Range("H122").Value = 4071.7827
$ws.Range("I122").Value = 3332.55
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 9997.650000000001
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -7547.650000000001
$ws.Range("N122").Value = -31900

$ws.Range("H132").Value = 5078.8696
$ws.Range("I132").Value = 3881.1333
$ws.Range("J132").Value = 7324.625
$ws.Range("K132").Value = 11643.3999
$ws.Range("L132").Value = 21973.875
$ws.Range("M132").Value = -9113.3999
$ws.Range("N132").Value = -27033.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H81").Value = 5333.3335
$ws.Range("J81").Value = 15000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -32122

$ws.Range("H84").Value = 5333.3335
$ws.Range("J84").Value = 15000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -160608

$ws.Range("H113").Value = 287.18182
$ws.Range("I113").Value = 320.25
$ws.Range("K113").Value = 960.75
$ws.Range("M113").Value = 1209.25
